$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the renamed Voting Classifier label first so the new shared-string
# pool keeps the same ordering as the source workbook (label right after
# "GaussianNB", ahead of the new confusion-matrix strings).
$ws.Range("B10").Value = "Voting Classifier -- Second Eval --  RF/GNB"

# Row 2 (Random Forest-Feature Selection-Tune=Recall)
$ws.Range("C2").Value = 0.7143
$ws.Range("D2").Value = 0.2353
$ws.Range("E2").Value = 0.3636
$ws.Range("F2").Value = 0.2857
$ws.Range("G2").Value = "[[46 13]`n [ 7  4]]"

# Row 3 (Gradient Boost-FeatureSelection-Tune=Recall)
$ws.Range("C3").Value = 0.7714
$ws.Range("D3").Value = 0.2222
$ws.Range("E3").Value = 0.1818
$ws.Range("F3").Value = 0.2
$ws.Range("G3").Value = "[[52  7]`n [ 9  2]]"

# Row 4 (Logistic Regression -- Second Eval - Untunned)
$ws.Range("C4").Value = 0.6429
$ws.Range("D4").Value = 0.15
$ws.Range("E4").Value = 0.2727
$ws.Range("F4").Value = 0.1935
$ws.Range("G4").Value = "[[42 17]`n [ 8  3]]"

# Row 5 (Random Forest -- Second Eval - Tune=Recall)
$ws.Range("C5").Value = 0.6714
$ws.Range("D5").Value = 0.2273
$ws.Range("E5").Value = 0.4545
$ws.Range("F5").Value = 0.303
$ws.Range("G5").Value = "[[42 17]`n [ 6  5]]"

# Row 6 (KNN -- Second Eval -- Tuning=Recall)
$ws.Range("C6").Value = 0.7429
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = "[[52  7]`n [11  0]]"

# Row 7 (LDA -- Second Eval - Tuning=Recall)
$ws.Range("C7").Value = 0.7286
$ws.Range("D7").Value = 0.25
$ws.Range("E7").Value = 0.3636
$ws.Range("F7").Value = 0.2963
$ws.Range("G7").Value = "[[47 12]`n [ 7  4]]"

# Row 8 (Classification Tree -- Second Eval)
$ws.Range("C8").Value = 0.6714
$ws.Range("D8").Value = 0.2
$ws.Range("E8").Value = 0.3636
$ws.Range("F8").Value = 0.2581
$ws.Range("G8").Value = "[[43 16]`n [ 7  4]]"

# Row 9 (GaussianNB)
$ws.Range("C9").Value = 0.4143
$ws.Range("D9").Value = 0.1875
$ws.Range("E9").Value = 0.8182
$ws.Range("F9").Value = 0.3051
$ws.Range("G9").Value = "[[20 39]`n [ 2  9]]"

# Row 10 (Voting Classifier -- now merged/renamed to "-- Second Eval --  RF/GNB")
$ws.Range("C10").Value = 0.4571
$ws.Range("D10").Value = 0.2
$ws.Range("E10").Value = 0.8182
$ws.Range("F10").Value = 0.3214
$ws.Range("G10").Value = "[[23 36]`n [ 2  9]]"

# Row 11 (old "Voting Classifier -- Second Eval -- Just CF/RF") is removed entirely
$ws.Rows.Item(11).Delete()
